$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 404.55173
$ws.Cells.Item(33, 10).Value = 1005.7143
$ws.Cells.Item(33, 12).Value = 1005.7143
$ws.Cells.Item(33, 14).Value = -1463.7143
$ws.Cells.Item(40, 8).Value = 3743.423
$ws.Cells.Item(40, 10).Value = 3739.75
$ws.Cells.Item(40, 12).Value = 3739.75
$ws.Cells.Item(40, 14).Value = -4089.75
$ws.Cells.Item(62, 8).Value = 97437680
$ws.Cells.Item(62, 9).Value = 105557070
$ws.Cells.Item(62, 11).Value = 105557070
$ws.Cells.Item(62, 13).Value = -105556446
$ws.Cells.Item(65, 8).Value = 97437680
$ws.Cells.Item(65, 9).Value = 105557070
$ws.Cells.Item(65, 11).Value = 527785350
$ws.Cells.Item(65, 13).Value = -527782230
$ws.Cells.Item(88, 8).Value = 2566.3635
$ws.Cells.Item(88, 10).Value = 2602.4285
$ws.Cells.Item(88, 12).Value = 2602.4285
$ws.Cells.Item(88, 14).Value = -3414.4285
$ws.Cells.Item(91, 8).Value = 2566.3635
$ws.Cells.Item(91, 10).Value = 2602.4285
$ws.Cells.Item(91, 12).Value = 2602.4285
$ws.Cells.Item(91, 14).Value = -5410.4285
$ws.Cells.Item(106, 8).Value = 300000
$ws.Cells.Item(106, 9).Value = 300000
$ws.Cells.Item(106, 11).Value = 300000
$ws.Cells.Item(106, 13).Value = -299369
$ws.Cells.Item(107, 8).Value = 22728830
$ws.Cells.Item(107, 9).Value = 27778486
$ws.Cells.Item(107, 10).Value = 5377.5
$ws.Cells.Item(107, 11).Value = 27778486
$ws.Cells.Item(107, 12).Value = 5377.5
$ws.Cells.Item(107, 13).Value = -27776566
$ws.Cells.Item(107, 14).Value = -9217.5
$ws.Cells.Item(141, 8).Value = 1166.9
$ws.Cells.Item(141, 9).Value = 925.82355
$ws.Cells.Item(141, 11).Value = 2777.47065
$ws.Cells.Item(141, 13).Value = 2402.52935

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 232.66667
$ws.Cells.Item(5, 10).Value = 67
$ws.Cells.Item(5, 12).Value = 67
$ws.Cells.Item(5, 14).Value = -291
$ws.Cells.Item(32, 8).Value = 4572490
$ws.Cells.Item(32, 9).Value = 4906525.5
$ws.Cells.Item(32, 11).Value = 4906525.5
$ws.Cells.Item(32, 13).Value = -4906238.5
$ws.Cells.Item(45, 8).Value = 6431.4707
$ws.Cells.Item(45, 9).Value = 8462.166999999999
$ws.Cells.Item(45, 11).Value = 8462.166999999999
$ws.Cells.Item(45, 13).Value = -8085.166999999999
$ws.Cells.Item(110, 8).Value = 11906631
$ws.Cells.Item(110, 9).Value = 19231580
$ws.Cells.Item(110, 11).Value = 19231580
$ws.Cells.Item(110, 13).Value = -19229535

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 232.66667
$ws.Cells.Item(4, 10).Value = 67
$ws.Cells.Item(4, 12).Value = 67
$ws.Cells.Item(4, 14).Value = -297
$ws.Cells.Item(86, 8).Value = 101361.65
$ws.Cells.Item(86, 9).Value = 1455
$ws.Cells.Item(86, 10).Value = 667499.3
$ws.Cells.Item(86, 11).Value = 1455
$ws.Cells.Item(86, 12).Value = 667499.3
$ws.Cells.Item(86, 13).Value = -332
$ws.Cells.Item(86, 14).Value = -669745.3
$ws.Cells.Item(89, 8).Value = 101361.65
$ws.Cells.Item(89, 9).Value = 1455
$ws.Cells.Item(89, 10).Value = 667499.3
$ws.Cells.Item(89, 11).Value = 7275
$ws.Cells.Item(89, 12).Value = 3337496.5
$ws.Cells.Item(89, 13).Value = -1659
$ws.Cells.Item(89, 14).Value = -3348728.5
$ws.Cells.Item(99, 8).Value = 1517.9445
$ws.Cells.Item(99, 10).Value = 2000.6666
$ws.Cells.Item(99, 12).Value = 2000.6666
$ws.Cells.Item(99, 14).Value = -4996.6666
$ws.Cells.Item(105, 8).Value = 45467948
$ws.Cells.Item(105, 9).Value = 58840124
$ws.Cells.Item(105, 10).Value = 2546
$ws.Cells.Item(105, 11).Value = 58840124
$ws.Cells.Item(105, 12).Value = 2546
$ws.Cells.Item(105, 13).Value = -58838377
$ws.Cells.Item(105, 14).Value = -6040

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1257.3334
$ws.Cells.Item(22, 9).Value = 805.5454999999999
$ws.Cells.Item(22, 10).Value = 2499.75
$ws.Cells.Item(22, 11).Value = 805.5454999999999
$ws.Cells.Item(22, 12).Value = 2499.75
$ws.Cells.Item(22, 13).Value = -455.5454999999999
$ws.Cells.Item(22, 14).Value = -3199.75
$ws.Cells.Item(62, 8).Value = 9092.333000000001
$ws.Cells.Item(62, 9).Value = 5104.6665
$ws.Cells.Item(62, 10).Value = 17067.666
$ws.Cells.Item(62, 11).Value = 5104.6665
$ws.Cells.Item(62, 12).Value = 17067.666
$ws.Cells.Item(62, 13).Value = -4480.6665
$ws.Cells.Item(62, 14).Value = -18315.666
$ws.Cells.Item(65, 8).Value = 9092.333000000001
$ws.Cells.Item(65, 9).Value = 5104.6665
$ws.Cells.Item(65, 10).Value = 17067.666
$ws.Cells.Item(65, 11).Value = 25523.3325
$ws.Cells.Item(65, 12).Value = 85338.33
$ws.Cells.Item(65, 13).Value = -22403.3325
$ws.Cells.Item(65, 14).Value = -91578.33
$ws.Cells.Item(134, 8).Value = 3290.561
$ws.Cells.Item(134, 9).Value = 2435.5715
$ws.Cells.Item(134, 11).Value = 7306.7145
$ws.Cells.Item(134, 13).Value = -4771.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 9).Value = 3258609
$ws.Cells.Item(4, 10).Value = 7075442
$ws.Cells.Item(4, 11).Value = 9775827
$ws.Cells.Item(4, 12).Value = 21226326
$ws.Cells.Item(4, 13).Value = -9775715
$ws.Cells.Item(4, 14).Value = -21226550
$ws.Cells.Item(5, 8).Value = 620.6923
$ws.Cells.Item(5, 10).Value = 1202.5
$ws.Cells.Item(5, 12).Value = 3607.5
$ws.Cells.Item(5, 14).Value = -3831.5
$ws.Cells.Item(99, 8).Value = 1281.125
$ws.Cells.Item(99, 9).Value = 1281.125
$ws.Cells.Item(99, 11).Value = 3843.375
$ws.Cells.Item(99, 13).Value = -1597.375
$ws.Cells.Item(135, 8).Value = 620.6923
$ws.Cells.Item(135, 10).Value = 1202.5
$ws.Cells.Item(135, 12).Value = 10822.5
$ws.Cells.Item(135, 14).Value = -15892.5
$ws.Cells.Item(137, 8).Value = 2457.5
$ws.Cells.Item(137, 10).Value = 2733.1428
$ws.Cells.Item(137, 12).Value = 8199.428400000001
$ws.Cells.Item(137, 14).Value = -18399.4284
$ws.Cells.Item(139, 8).Value = 83335970
$ws.Cells.Item(139, 9).Value = 125001700
$ws.Cells.Item(139, 10).Value = 4500
$ws.Cells.Item(139, 11).Value = 375005100
$ws.Cells.Item(139, 12).Value = 13500
$ws.Cells.Item(139, 13).Value = -374999960
$ws.Cells.Item(139, 14).Value = -23780

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 12730.462
$ws.Cells.Item(113, 10).Value = 1010
$ws.Cells.Item(113, 12).Value = 1010
$ws.Cells.Item(113, 14).Value = -5350
$ws.Cells.Item(122, 8).Value = 10262.722
$ws.Cells.Item(122, 10).Value = 12994.728
$ws.Cells.Item(122, 12).Value = 38984.18399999999
$ws.Cells.Item(122, 14).Value = -43884.18399999999
$ws.Cells.Item(126, 8).Value = 2722.92
$ws.Cells.Item(126, 9).Value = 2142.6667
$ws.Cells.Item(126, 11).Value = 6428.000100000001
$ws.Cells.Item(126, 13).Value = -3958.000100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3239.1794
$ws.Cells.Item(22, 9).Value = 2018.8636
$ws.Cells.Item(22, 10).Value = 4818.4116
$ws.Cells.Item(22, 11).Value = 2018.8636
$ws.Cells.Item(22, 12).Value = 4818.4116
$ws.Cells.Item(22, 13).Value = -1723.8636
$ws.Cells.Item(22, 14).Value = -5408.4116
$ws.Cells.Item(27, 8).Value = 3239.1794
$ws.Cells.Item(27, 9).Value = 2018.8636
$ws.Cells.Item(27, 10).Value = 4818.4116
$ws.Cells.Item(27, 11).Value = 2018.8636
$ws.Cells.Item(27, 12).Value = 4818.4116
$ws.Cells.Item(27, 13).Value = -1911.8636
$ws.Cells.Item(27, 14).Value = -5032.4116
$ws.Cells.Item(100, 8).Value = 22729504
$ws.Cells.Item(100, 9).Value = 250000000
$ws.Cells.Item(100, 10).Value = 2453.5
$ws.Cells.Item(100, 11).Value = 250000000
$ws.Cells.Item(100, 12).Value = 2453.5
$ws.Cells.Item(100, 13).Value = -249999459
$ws.Cells.Item(100, 14).Value = -3535.5
$ws.Cells.Item(132, 8).Value = 5770.5835
$ws.Cells.Item(132, 9).Value = 4562.9473
$ws.Cells.Item(132, 10).Value = 10359.6
$ws.Cells.Item(132, 11).Value = 13688.8419
$ws.Cells.Item(132, 12).Value = 31078.8
$ws.Cells.Item(132, 13).Value = -11158.8419
$ws.Cells.Item(132, 14).Value = -36138.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2299.8
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 14).Value = -13900
$ws.Cells.Item(132, 8).Value = 3686
$ws.Cells.Item(132, 9).Value = 2183.3408
$ws.Cells.Item(132, 11).Value = 6550.0224
$ws.Cells.Item(132, 13).Value = -4020.0224
